# Consolidate findings in v2.1.0_spec
#
# Slide 2 ("EmbeddingStatus" summary block, under the "Update Embedding
# Status" arrow):
#  - widen the "ApplicationName / ReleaseNumber / EmbeddingStatus" label
#    textbox so there is room to also label a "ReasonOfFailure" value
#  - tidy the "EmbeddingStatus" run (drop its trailing space) and extend
#    the trailing ":" run into ":      ReasonOfFailure:"
#  - add a new empty rectangle (styled like the neighboring status value
#    boxes) to the right, to hold the ReasonOfFailure value

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- widen the label textbox ("Textfeld 5") ---------------------------
$labelBox = $s.Shapes.Item("Textfeld 5")
$labelBox.Width = 620.5794

# --- edit the text in-place, preserving the existing run formatting ---
$tr = $labelBox.TextFrame.TextRange
$full = $tr.Text

$oldEmbed = "EmbeddingStatus "
$idx = $full.IndexOf($oldEmbed)
$embedRun = $tr.Characters($idx + 1, $oldEmbed.Length)
$embedRun.Text = "EmbeddingStatus"

$full = $tr.Text
$idxColon = $full.LastIndexOf(":")
$colonRun = $tr.Characters($idxColon + 1, 1)
$colonRun.Text = ":      ReasonOfFailure:"

# --- add the new "ReasonOfFailure" value rectangle ---------------------
# Duplicate one of the existing small status-value rectangles so the new
# shape inherits the same style (line/fill/effect/font refs) and empty,
# center-aligned text body, then move/resize and rename it.
$template = $s.Shapes.Item("Rechteck 7")
$newShapeRange = $template.Duplicate()
$newShape = $newShapeRange.Item(1)
$newShape.Name = "Rechteck 21"
$newShape.Left = 562.3709
$newShape.Top = 455.8374
$newShape.Width = 195.912
$newShape.Height = 26.1981
